$d = $word.ActiveDocument

# Locate the paragraph that ends the intro ("Markdown is just a plain text
# format...") so we can insert the new sentence right after it, before the
# "For example:" paragraph.
$anchorText = "Markdown is just a plain text format that is designed to be " +
    "easy to write, and, even more importantly, easy to read."

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq $anchorText) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the target paragraph to anchor the insertion."
}

$anchorEnd = $target.Range.End

# Insert a brand-new paragraph mark right after the target paragraph.
$anchor = $target.Range
$anchor.Collapse(0)          # wdCollapseEnd
$anchor.InsertParagraphAfter()

# Find the paragraph that was just created: the one that now starts where
# the target paragraph used to end.
$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $anchorEnd) {
        $newPara = $p
        break
    }
}

$newPara.Range.InsertBefore("It just requires a little training but then it is very fast to edit.")
$newPara.Style = "Body Text"
